$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
$src = $ws.Range("A10:B10")
$dst = $ws.Range("A11:B11")
$src.Copy()
$dst.PasteSpecial(-4122)
